$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at I (shifts old I=date, J=legislator_name, K=legislator_id
# one column to the right, becoming J, K, L respectively); formatting carries
# over automatically for the inserted column.
$ws.Range("I1:I4").EntireColumn.Insert()

# Fill the new "category" column (I)
$ws.Range("I1").Value = "category"
$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"
$ws.Range("I4").Value = "normal"

# Add the two new trailing columns: source_file (M) and index (N).
# Copy formatting from the existing header/data columns (H header style, H2:H4 data style)
# before writing the new values.
$ws.Range("H1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)
$ws.Range("H2:H4").Copy()
$ws.Range("M2:M4").PasteSpecial(-4122)
$ws.Range("N2:N4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("M1").Value = "source_file"
$ws.Range("M2").Value = "tmp56941"
$ws.Range("M3").Value = "tmp56941"
$ws.Range("M4").Value = "tmp56941"

$ws.Range("N1").Value = "index"
$ws.Range("N2").Value = 89
$ws.Range("N3").Value = 90
$ws.Range("N4").Value = 91
